$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Continue the daily-count table with the next date column (BW = "22-sep"),
# matching the formatting already used for the preceding date columns.
$ws.Range("BW1").NumberFormat = "@"
$ws.Range("BW1").Value = "22-sep"

$ws.Range("BW2:BW11").HorizontalAlignment = -4108
$ws.Range("BW2:BW11").NumberFormat = "0"

$ws.Range("BW2").Value = 12
$ws.Range("BW3").Value = 13
$ws.Range("BW4").Value = 12
$ws.Range("BW5").Value = 11
$ws.Range("BW6").Value = 8
$ws.Range("BW7").Value = 14
$ws.Range("BW8").Value = 19
$ws.Range("BW9").Value = 14
$ws.Range("BW10").Value = 15
$ws.Range("BW11").Value = 3

# Match the final cell selection recorded in the saved file.
$ws.Range("BY5").Select()
